$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiaJourney")

$ws.Cells.Item(4, 2).Value2 = 'watch me let go... this is only for you 💕'
$ws.Cells.Item(6, 2).Value2 = 'cum with me... I''m letting go right now, watch me'
$ws.Cells.Item(7, 2).Value2 = 'I''m so close... I can feel it building everywhere and I don''t want to cum alone 💕'
$ws.Cells.Item(8, 2).Value2 = 'oh god I can''t hold on 🌸'
$ws.Cells.Item(9, 2).Value2 = 'this is what you made me do and you need to see every second of it 💕'
$ws.Cells.Item(11, 2).Value2 = 'I need you to watch what I''m doing right now... you have to see this'
$ws.Cells.Item(12, 2).Value2 = 'my fingers keep going deeper and faster and my whole body is shaking 💕'
$ws.Cells.Item(13, 2).Value2 = 'I''m touching my pussy and imagining it''s your hands on me... I need more'
$ws.Cells.Item(14, 2).Value2 = 'fuck 🌸'
$ws.Cells.Item(15, 2).Value2 = 'look what you did to me... I couldn''t stop 💕'
$ws.Cells.Item(17, 2).Value2 = 'tell me what you want me to do to myself right now... I''ll do anything you say'
$ws.Cells.Item(18, 2).Value2 = 'I''m so wet right now babe... you have no idea what your words do to my body 💕'
$ws.Cells.Item(19, 2).Value2 = 'but I can''t stop now... my fingers are between my legs and it''s all because of you'
$ws.Cells.Item(20, 2).Value2 = 'oh god... I can''t believe I just did that 🌸'
$ws.Cells.Item(21, 2).Value2 = 'I want you to see what you''re doing to me right now 💕'
$ws.Cells.Item(23, 2).Value2 = 'my hand keeps sliding lower and I can''t stop it babe... I don''t even want to'
$ws.Cells.Item(24, 2).Value2 = 'I keep running my hands down my body and everything is so sensitive... it''s like every touch is amplified because of you'
$ws.Cells.Item(25, 2).Value2 = 'you really liked that? knowing you saw me is making my heart race so fast right now 🌸'
